$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$v = $ws1.Range("B2").Value()
Write-Output "VALUE=$v"
